# Hortaliza, Vega Monumental Concepción - Apio
# Weekly update: insert two new price records (week of 44785) at the top of
# the data block (rows 258-259), pushing the existing rows down by two and
# extending the used range to A1:R321.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 258.
$ws.Rows("258:259").Insert()

# New row 258 - "Primera" quality record for the new date.
$ws.Range("A258").Value = 11
$ws.Range("B258").Value = "Vega Monumental Concepción"
$ws.Range("C258").Value = "Bíobío"
$ws.Range("D258").Value = 44785
$ws.Range("E258").Value = 8
$ws.Range("F258").Value = 100112017
$ws.Range("G258").Value = "Apio"
$ws.Range("H258").Value = "Americana (o)"
$ws.Range("I258").Value = "Primera"
$ws.Range("J258").Value = 100
$ws.Range("K258").Value = 9000
$ws.Range("L258").Value = 9500
$ws.Range("M258").Value = 9250
$ws.Range("N258").Value = "$/docena de matas"
$ws.Range("O258").Value = "Región de Coquimbo"
$ws.Range("P258").Value = 1542
$ws.Range("Q258").Value = 6
$ws.Range("R258").Value = "Hortaliza"

# New row 259 - "Segunda" quality record for the new date.
$ws.Range("A259").Value = 11
$ws.Range("B259").Value = "Vega Monumental Concepción"
$ws.Range("C259").Value = "Bíobío"
$ws.Range("D259").Value = 44785
$ws.Range("E259").Value = 8
$ws.Range("F259").Value = 100112017
$ws.Range("G259").Value = "Apio"
$ws.Range("H259").Value = "Americana (o)"
$ws.Range("I259").Value = "Segunda"
$ws.Range("J259").Value = 50
$ws.Range("K259").Value = 8000
$ws.Range("L259").Value = 8000
$ws.Range("M259").Value = 8000
$ws.Range("N259").Value = "$/docena de matas"
$ws.Range("O259").Value = "Región de Coquimbo"
$ws.Range("P259").Value = 1333
$ws.Range("Q259").Value = 6
$ws.Range("R259").Value = "Hortaliza"
